# Weekly update: insert a new week's price record at the top of the
# "Poroto verde" (Feria Lagunitas de Puerto Montt) data block, pushing the
# existing rows (59..81) down by one (to 60..82) and growing the sheet's
# used range from R81 to R82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 59-81 down to 60-82, inserting a blank row 59.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with this week's record.
$ws.Cells.Item(59, 1).Value  = 4
$ws.Cells.Item(59, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value  = "Los Lagos"
$ws.Cells.Item(59, 4).Value  = 44719
$ws.Cells.Item(59, 5).Value  = 10
$ws.Cells.Item(59, 6).Value  = 100112031
$ws.Cells.Item(59, 7).Value  = "Poroto verde"
$ws.Cells.Item(59, 8).Value  = "Magnum"
$ws.Cells.Item(59, 9).Value  = "Primera"
$ws.Cells.Item(59, 10).Value = 35
$ws.Cells.Item(59, 11).Value = 35000
$ws.Cells.Item(59, 12).Value = 35000
$ws.Cells.Item(59, 13).Value = 35000
$ws.Cells.Item(59, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(59, 15).Value = "Perú"
$ws.Cells.Item(59, 16).Value = 1400
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"
